$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.610.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.363.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "661.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -9.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.422"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -10.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.358.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.210"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.372.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000257"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.972.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.367.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.560"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +28.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "503.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000201"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.542.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.149"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.995"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +19.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.555"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "527.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  -5.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.849"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0424"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.20%  "
